$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2458.2454
$ws.Range("I15").Value = 2458.2454
$ws.Range("K15").Value = 7374.736199999999
$ws.Range("M15").Value = -7205.736199999999

$ws.Range("H42").Value = 190.41667
$ws.Range("I42").Value = 190.41667
$ws.Range("K42").Value = 571.25001
$ws.Range("M42").Value = -341.25001

$ws.Range("H43").Value = 151975.36
$ws.Range("J43").Value = 191236.9
$ws.Range("L43").Value = 191236.9
$ws.Range("N43").Value = -191374.9

$ws.Range("H69").Value = 53753.75
$ws.Range("I69").Value = 30000
$ws.Range("K69").Value = 90000
$ws.Range("M69").Value = -89126

$ws.Range("H72").Value = 53753.75
$ws.Range("I72").Value = 30000
$ws.Range("K72").Value = 270000
$ws.Range("M72").Value = -265632

$ws.Range("H76").Value = 3891.111
$ws.Range("I76").Value = 3956.5
$ws.Range("K76").Value = 3956.5
$ws.Range("M76").Value = -3641.5

$ws.Range("H79").Value = 3891.111
$ws.Range("I79").Value = 3956.5
$ws.Range("K79").Value = 3956.5
$ws.Range("M79").Value = -2864.5

$ws.Range("H101").Value = 1344.4166
$ws.Range("I101").Value = 1266.75
$ws.Range("K101").Value = 3800.25
$ws.Range("M101").Value = -2178.25

$ws.Range("H137").Value = 4271.5483
$ws.Range("I137").Value = 2766.5881
$ws.Range("J137").Value = 6099
$ws.Range("K137").Value = 8299.764299999999
$ws.Range("L137").Value = 18297
$ws.Range("M137").Value = -5749.764299999999
$ws.Range("N137").Value = -23397

$ws.Range("H138").Value = 3633.5278
$ws.Range("I138").Value = 1736
$ws.Range("K138").Value = 5208
$ws.Range("M138").Value = -68

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8198802
$ws.Range("I32").Value = 9092218
$ws.Range("K32").Value = 9092218
$ws.Range("M32").Value = -9091931

$ws.Range("H64").Value = 49500
$ws.Range("J64").Value = 49500
$ws.Range("L64").Value = 49500
$ws.Range("N64").Value = -49996

$ws.Range("H67").Value = 49500
$ws.Range("J67").Value = 49500
$ws.Range("L67").Value = 49500
$ws.Range("N67").Value = -51216

$ws.Range("H74").Value = 7227415.5
$ws.Range("I74").Value = 9261873
$ws.Range("K74").Value = 9261873
$ws.Range("M74").Value = -9260999

$ws.Range("H77").Value = 7227415.5
$ws.Range("I77").Value = 9261873
$ws.Range("K77").Value = 46309365
$ws.Range("M77").Value = -46304997

$ws.Range("H97").Value = 1893.8572
$ws.Range("I97").Value = 1893.8572
$ws.Range("K97").Value = 1893.8572
$ws.Range("M97").Value = -1397.8572

$ws.Range("H132").Value = 3203.8293
$ws.Range("I132").Value = 1483.4849
$ws.Range("K132").Value = 4450.4547
$ws.Range("M132").Value = -1920.4547

$ws.Range("H138").Value = 15000
$ws.Range("I138").Value = 15000
$ws.Range("K138").Value = 15000
$ws.Range("M138").Value = -9860

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 548.75
$ws.Range("I94").Value = 482.0909
$ws.Range("J94").Value = 695.4
$ws.Range("K94").Value = 482.0909
$ws.Range("L94").Value = 695.4
$ws.Range("M94").Value = -31.09089999999998
$ws.Range("N94").Value = -1597.4

$ws.Range("H134").Value = 371810.97
$ws.Range("I134").Value = 1355.28
$ws.Range("J134").Value = 5002507
$ws.Range("K134").Value = 4065.84
$ws.Range("L134").Value = 15007521
$ws.Range("M134").Value = -1530.84
$ws.Range("N134").Value = -15012591

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 865
$ws.Range("I16").Value = 724.25
$ws.Range("K16").Value = 724.25
$ws.Range("M16").Value = -437.25

$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()

$ws.Range("H86").Value = 3899.6
$ws.Range("I86").Value = 3899.6
$ws.Range("K86").Value = 3899.6
$ws.Range("M86").Value = -2776.6

$ws.Range("H89").Value = 3899.6
$ws.Range("I89").Value = 3899.6
$ws.Range("K89").Value = 19498
$ws.Range("M89").Value = -13882

$ws.Range("H107").Value = 955.46155
$ws.Range("I107").Value = 527.2
$ws.Range("J107").Value = 1223.125
$ws.Range("K107").Value = 527.2
$ws.Range("L107").Value = 1223.125
$ws.Range("M107").Value = 1392.8
$ws.Range("N107").Value = -5063.125

$ws.Range("H113").Value = 865
$ws.Range("I113").Value = 724.25
$ws.Range("K113").Value = 724.25
$ws.Range("M113").Value = 1445.75

$ws.Range("H122").Value = 1716.125
$ws.Range("I122").Value = 1758.6666
$ws.Range("K122").Value = 5275.9998
$ws.Range("M122").Value = -2825.9998

$ws.Range("H132").Value = 2483.6316
$ws.Range("I132").Value = 1952.3529
$ws.Range("J132").Value = 6999.5
$ws.Range("K132").Value = 5857.0587
$ws.Range("L132").Value = 20998.5
$ws.Range("M132").Value = -3327.0587
$ws.Range("N132").Value = -26058.5

$ws.Range("H134").Value = 335009.7
$ws.Range("I134").Value = 358796.1
$ws.Range("K134").Value = 1076388.3
$ws.Range("M134").Value = -1073853.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 2000
$ws.Range("I57").Value = 2000
$ws.Range("K57").Value = 6000
$ws.Range("M57").Value = -5441

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 92488.5
$ws.Range("I5").Value = 109983
$ws.Range("J5").Value = 40005
$ws.Range("K5").Value = 109983
$ws.Range("L5").Value = 40005
$ws.Range("M5").Value = -109871
$ws.Range("N5").Value = -40229

$ws.Range("H39").Value = 75000
$ws.Range("J39").Value = 75000
$ws.Range("L39").Value = 75000
$ws.Range("N39").Value = -76064

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 50001484
$ws.Range("I93").Value = 62501580
$ws.Range("K93").Value = 62501580
$ws.Range("M93").Value = -62500332

$ws.Range("H119").Value = 82330
$ws.Range("J119").Value = 82330
$ws.Range("L119").Value = 82330
$ws.Range("N119").Value = -92006

$ws.Range("H122").Value = 5330.1797
$ws.Range("I122").Value = 4980.8438
$ws.Range("K122").Value = 14942.5314
$ws.Range("M122").Value = -12492.5314

$ws.Range("H129").Value = 83195
$ws.Range("I129").Value = 66390
$ws.Range("J129").Value = 100000
$ws.Range("K129").Value = 66390
$ws.Range("L129").Value = 100000
$ws.Range("M129").Value = -61390
$ws.Range("N129").Value = -110000

$ws.Range("H132").Value = 174704.06
$ws.Range("I132").Value = 1843.3334
$ws.Range("J132").Value = 1004435.6
$ws.Range("K132").Value = 5530.0002
$ws.Range("L132").Value = 3013306.8
$ws.Range("M132").Value = -3000.0002
$ws.Range("N132").Value = -3018366.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 4500
$ws.Range("I55").Value = 3000
$ws.Range("J55").Value = 6000
$ws.Range("K55").Value = 3000
$ws.Range("L55").Value = 6000
$ws.Range("M55").Value = -2723
$ws.Range("N55").Value = -6554

$ws.Range("H126").Value = 1549.8422
$ws.Range("I126").Value = 1585.4375
$ws.Range("J126").Value = 1360
$ws.Range("K126").Value = 4756.3125
$ws.Range("L126").Value = 4080
$ws.Range("M126").Value = -2286.3125
$ws.Range("N126").Value = -9020

$ws.Range("H129").Value = 124000
$ws.Range("J129").Value = 124000
$ws.Range("L129").Value = 124000
$ws.Range("N129").Value = -134000

$ws.Range("H132").Value = 2464.5334
$ws.Range("I132").Value = 2053.1853
$ws.Range("K132").Value = 6159.5559
$ws.Range("M132").Value = -3629.5559
